$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Obrigatorio" column (E) from "N" to "S" for rows 2-9
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = "S"
}

# Add new row 10: NFE12-FILLER
$ws.Cells.Item(10, 1).Value = "NFE12-FILLER"
$ws.Cells.Item(10, 2).Value = 47
$ws.Cells.Item(10, 3).Value = 494
$ws.Cells.Item(10, 4).Value = "TEXTO"
$ws.Cells.Item(10, 5).Value = "N"
$ws.Cells.Item(10, 6).Value = ""
